# Update "想去人数" (want-to-go count) values on both the "展览" and
# "全部类型" sheets. These two sheets contain identical data tables,
# so the same four cells (F3, F8, F9, F11) are updated on each.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 797
    $ws.Range("F8").Value = 2031
    $ws.Range("F9").Value = 7155
    $ws.Range("F11").Value = 409
}
